# ============================================================
# Actualizacion desde MV -datos-
# Adds two new daily FX rows (04-10-2021, 05-10-2021) and fixes
# a previously-missing value in D190.
# ============================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix D190: was 0, should be 196.67 ---
$ws.Cells.Item(190, 4).Value = 196.67

# --- Row 191 ---
$ws.Cells.Item(191, 1).Formula = "=""04-10-2021"""
$ws.Cells.Item(191, 1).Copy()
$ws.Cells.Item(191, 1).PasteSpecial(-4163)
$ws.Cells.Item(191, 2).Value = 23.91
$ws.Cells.Item(191, 3).Value = 803.9
$ws.Cells.Item(191, 4).Value = 194.73
$ws.Cells.Item(191, 5).Value = 117.19
$ws.Cells.Item(191, 6).Value = 1.29
$ws.Cells.Item(191, 7).Value = 36.84
$ws.Cells.Item(191, 8).Value = 125.38
$ws.Cells.Item(191, 9).Value = 6.21
$ws.Cells.Item(191, 10).Value = 93.25
$ws.Cells.Item(191, 11).Value = 91.90000000000001
$ws.Cells.Item(191, 12).Value = 1134.81
$ws.Cells.Item(191, 13).Value = 218.93
$ws.Cells.Item(191, 14).Value = 88.78
$ws.Cells.Item(191, 15).Value = 584.74
$ws.Cells.Item(191, 16).Value = 636
$ws.Cells.Item(191, 17).Value = 803.9
$ws.Cells.Item(191, 18).Value = 803.9
$ws.Cells.Item(191, 19).Value = 380.56
$ws.Cells.Item(191, 20).Value = 980.37
$ws.Cells.Item(191, 21).Value = 803.9
$ws.Cells.Item(191, 22).Value = 593.33
$ws.Cells.Item(191, 23).Value = 103.26
$ws.Cells.Item(191, 24).Value = 558.23
$ws.Cells.Item(191, 25).Value = 28.88
$ws.Cells.Item(191, 26).Value = 0.04
$ws.Cells.Item(191, 27).Value = 932.38
$ws.Cells.Item(191, 28).Value = 2.61
$ws.Cells.Item(191, 29).Value = 7.84
$ws.Cells.Item(191, 30).Value = 864.6900000000001
$ws.Cells.Item(191, 31).Value = 0.12
$ws.Cells.Item(191, 32).Value = 30.21
$ws.Cells.Item(191, 33).Value = 188.63
$ws.Cells.Item(191, 34).Value = 51.24
$ws.Cells.Item(191, 35).Value = 1090.18
$ws.Cells.Item(191, 36).Value = 90.81
$ws.Cells.Item(191, 37).Value = 194.76
$ws.Cells.Item(191, 38).Value = 8.140000000000001
$ws.Cells.Item(191, 39).Value = 0.21
$ws.Cells.Item(191, 40).Value = 33.5
$ws.Cells.Item(191, 41).Value = 14.3
$ws.Cells.Item(191, 42).Value = 15.83
$ws.Cells.Item(191, 43).Value = 39.35
$ws.Cells.Item(191, 44).Value = 18.79
$ws.Cells.Item(191, 45).Value = 104.02
$ws.Cells.Item(191, 46).Value = 54.11
$ws.Cells.Item(191, 47).Value = 149.84
$ws.Cells.Item(191, 48).Value = 0.02
$ws.Cells.Item(191, 49).Value = 214.37
$ws.Cells.Item(191, 50).Value = 192.46
$ws.Cells.Item(191, 51).Value = 11.06
$ws.Cells.Item(191, 52).Value = 10.85
$ws.Cells.Item(191, 53).Value = 0.06
$ws.Cells.Item(191, 54).Value = 4.74
$ws.Cells.Item(191, 55).Value = 250.05
$ws.Cells.Item(191, 56).Value = 1.89
$ws.Cells.Item(191, 57).Value = 0.68
$ws.Cells.Item(191, 58).Value = 7.24
$ws.Cells.Item(191, 59).Value = 125.02
$ws.Cells.Item(191, 60).Value = 203.68

# --- Row 192 ---
$ws.Cells.Item(192, 1).Formula = "=""05-10-2021"""
$ws.Cells.Item(192, 1).Copy()
$ws.Cells.Item(192, 1).PasteSpecial(-4163)
$ws.Cells.Item(192, 2).Value = 23.88
$ws.Cells.Item(192, 3).Value = 805.89
$ws.Cells.Item(192, 4).Value = 193.2
$ws.Cells.Item(192, 5).Value = 117.48
$ws.Cells.Item(192, 6).Value = 1.29
$ws.Cells.Item(192, 7).Value = 36.93
$ws.Cells.Item(192, 8).Value = 125.9
$ws.Cells.Item(192, 9).Value = 6.3
$ws.Cells.Item(192, 10).Value = 93.90000000000001
$ws.Cells.Item(192, 11).Value = 92.34
$ws.Cells.Item(192, 12).Value = 1138.91
$ws.Cells.Item(192, 13).Value = 219.47
$ws.Cells.Item(192, 14).Value = 89.52
$ws.Cells.Item(192, 15).Value = 587
$ws.Cells.Item(192, 16).Value = 640.61
$ws.Cells.Item(192, 17).Value = 805.89
$ws.Cells.Item(192, 18).Value = 805.89
$ws.Cells.Item(192, 19).Value = 382.48
$ws.Cells.Item(192, 20).Value = 982.79
$ws.Cells.Item(192, 21).Value = 805.89
$ws.Cells.Item(192, 22).Value = 593.83
$ws.Cells.Item(192, 23).Value = 103.49
$ws.Cells.Item(192, 24).Value = 561.13
$ws.Cells.Item(192, 25).Value = 28.92
$ws.Cells.Item(192, 26).Value = 0.04
$ws.Cells.Item(192, 27).Value = 936.3200000000001
$ws.Cells.Item(192, 28).Value = 2.63
$ws.Cells.Item(192, 29).Value = 7.87
$ws.Cells.Item(192, 30).Value = 871.33
$ws.Cells.Item(192, 31).Value = 0.12
$ws.Cells.Item(192, 32).Value = 30.57
$ws.Cells.Item(192, 33).Value = 189.25
$ws.Cells.Item(192, 34).Value = 51.46
$ws.Cells.Item(192, 35).Value = 1097.2
$ws.Cells.Item(192, 36).Value = 90.95
$ws.Cells.Item(192, 37).Value = 195.13
$ws.Cells.Item(192, 38).Value = 8.15
$ws.Cells.Item(192, 39).Value = 0.21
$ws.Cells.Item(192, 40).Value = 33.58
$ws.Cells.Item(192, 41).Value = 14.34
$ws.Cells.Item(192, 42).Value = 15.9
$ws.Cells.Item(192, 43).Value = 39.21
$ws.Cells.Item(192, 44).Value = 18.81
$ws.Cells.Item(192, 45).Value = 104.27
$ws.Cells.Item(192, 46).Value = 53.51
$ws.Cells.Item(192, 47).Value = 148.26
$ws.Cells.Item(192, 48).Value = 0.02
$ws.Cells.Item(192, 49).Value = 214.9
$ws.Cells.Item(192, 50).Value = 193.12
$ws.Cells.Item(192, 51).Value = 11.12
$ws.Cells.Item(192, 52).Value = 10.85
$ws.Cells.Item(192, 53).Value = 0.06
$ws.Cells.Item(192, 54).Value = 4.73
$ws.Cells.Item(192, 55).Value = 250
$ws.Cells.Item(192, 56).Value = 1.9
$ws.Cells.Item(192, 57).Value = 0.68
$ws.Cells.Item(192, 58).Value = 7.27
$ws.Cells.Item(192, 59).Value = 124.96
$ws.Cells.Item(192, 60).Value = 203.85

$excel.CutCopyMode = 0
